# Logout Test case update: add new "ExptectedResult" column (F) to the Logout sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logout")

# New header cell, formatted like the other header cells (bold / filled)
$ws.Range("F1").Value = "ExptectedResult"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data cells for the two existing rows
$ws.Range("F2").Value = "Logged out from user toggle button."
$ws.Range("F3").Value = "Logged out from Account web page."

# Size the new column to fit its contents
$ws.Columns.Item(6).ColumnWidth = 30

# Match the active selection left behind in the source file
$ws.Range("F3").Select()
